# Weekly price update: insert a new weekly record at the top of the data
# (row 55), pushing all subsequent rows down by one. This mirrors the
# author's commit "Fruta / hortaliza, semanal" which adds the latest
# week's Albahaca price observation for Terminal La Palmera de La Serena.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 55; existing rows 55-107 shift to 56-108.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row with the latest week's observation.
$ws.Cells.Item(55, 1).Value2 = 8
$ws.Cells.Item(55, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(55, 3).Value2 = "Coquimbo"
$ws.Cells.Item(55, 4).Value2 = 44778
$ws.Cells.Item(55, 5).Value2 = 4
$ws.Cells.Item(55, 6).Value2 = 100112052
$ws.Cells.Item(55, 7).Value2 = "Albahaca"
$ws.Cells.Item(55, 8).Value2 = "Sin especificar"
$ws.Cells.Item(55, 9).Value2 = "Primera"
$ws.Cells.Item(55, 10).Value2 = 1120
$ws.Cells.Item(55, 11).Value2 = 3500
$ws.Cells.Item(55, 12).Value2 = 4000
$ws.Cells.Item(55, 13).Value2 = 3750
$ws.Cells.Item(55, 14).Value2 = "`$/paquete"
$ws.Cells.Item(55, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(55, 16).Value2 = 3750
$ws.Cells.Item(55, 17).Value2 = 1
$ws.Cells.Item(55, 18).Value2 = "Hortaliza"
